# Regenerate merged AHB files
# - Rename header row: "<Name>_old" -> "<Name>_FV2210", "<Name>_new" -> "<Name>_FV2304"
# - Wrap the used range in a native Excel Table ("Table1")
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A..J (1..10) carry the "_old" -> "_FV2210" header names.
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $baseNames[$i] + "_FV2210"
}

# Column K (11) is the unchanged "diff" header.

# Columns L..U (12..21) carry the "_new" -> "_FV2304" header names.
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = $baseNames[$i] + "_FV2304"
}

# Turn the used range A1:U57 into a native table ("Table1").
$tableRange = $ws.Range("A1:U57")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# Freeze the header row.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
